# Generate and process indicate message
# Moves existing "OTA update..." block content down by two rows (53/54 -> 56/57,
# 55 -> 57, 57 -> 59, 58 -> 60, 60 -> 62, 62 -> 64, 64 -> 66, 66 -> 68, 68 -> 70,
# 70 -> 72, 72 -> 74, 74 -> 76) to make room for a new "Indicate (identify node)"
# control command block inserted at rows 53-54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-Range([string]$srcRange, [string]$dstRange) {
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4163)  # xlPasteValues
    $ws.Range($srcRange).Clear()
}

# Work bottom-up so a source row is always read before it is overwritten.
Move-Range "D74:G74" "D76:G76"
Move-Range "D72:G72" "D74:G74"
Move-Range "D70:G70" "D72:G72"
Move-Range "D68:G68" "D70:G70"
Move-Range "D66:G66" "D68:G68"
Move-Range "D64:G64" "D66:G66"
Move-Range "D62:N62" "D64:N64"
Move-Range "C60:F60" "C62:F62"
Move-Range "D58:N58" "D60:N60"
Move-Range "E57:E57" "E59:E59"
Move-Range "G55:L55" "G57:L57"
Move-Range "H54:H54" "H56:H56"

$ws.Application.CutCopyMode = 0

# New row 53: label for the "Indicate" control command.
$ws.Range("F53").Style = "Normal"
$ws.Range("F53").Value = $null
$ws.Range("H53").Value = "Indicate (identify node)"

# New row 54: direction arrow + "04" code, matching the existing
# "Get version" / "Get sleep time" / "Set sleep time" block styling.
$ws.Range("F54").Style = "Normal"
$ws.Range("F54").Value = $null
$ws.Range("G47").Copy()
$ws.Range("G54").PasteSpecial(-4122)
$ws.Range("G54").Value = "<------------------------------------"
$ws.Range("H47").Copy()
$ws.Range("H54").PasteSpecial(-4122)
$ws.Range("H54").Value = "04"

$ws.Application.CutCopyMode = 0

# Fix up the header area that moved: row 54 used to hold the "OTA update"
# label which is now two rows further down at row 56; row 54 instead shows
# the new Indicate block, so clear the stale ignoredErrors / selection state
# automatically handled by recalculation. Update sheet selection + dimension.
$ws.Range("J54").Select()

Write-Output "done"
